$d = $word.ActiveDocument

# --- Step 1: insert a brand-new paragraph before "Bei einer linearen Regression ..."
#     containing the Gaussian-error-propagation note, and move the _GoBack
#     bookmark onto it.
$p6 = $d.Paragraphs(6).Range
$p6.InsertParagraphBefore()
$newPara = $d.Paragraphs(6).Range
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Gau&#223;sche</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> Fehlerfortpflanzung ist handschriftlich hinzugef&#252;gt.</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$newPara.InsertXML($newParaXml)

# --- Step 2: split "ax+b" into its own (spell-checked) run inside the
#     "Bei einer linearen Regression ..." paragraph (now paragraph 7).
$p7 = $d.Paragraphs(7).Range
$p7Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r><w:t>Bei einer linearen Regression y=</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>ax+b</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> werden die Parameter a und b nach </w:t></w:r>' + `
  '</w:p>'
$p7.InsertXML($p7Xml)

# --- Step 3: split "IPython" into its own (spell-checked) run inside the
#     "Die lineare Regression ..." paragraph (now paragraph 11), and drop
#     the bookmark that used to sit at the end of this paragraph (it now
#     lives on the new paragraph inserted in step 1).
$p11 = $d.Paragraphs(11).Range
$p11Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r><w:t xml:space="preserve">Die lineare Regression wurde mit </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>IPython</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> 5.3</w:t></w:r>' + `
  '<w:r><w:t>.0 in Python 3.6.1 durchgef&#252;hrt, wobei die Unsicherheiten f&#252;r a und b auch bestimmt werden.</w:t></w:r>' + `
  '</w:p>'
$p11.InsertXML($p11Xml)
